$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.34
$ws.Range("E3").Value = 1.3
$ws.Range("E4").Value = 1.22
$ws.Range("C5").Value = 1.35
$ws.Range("D5").Value = 1.35
$ws.Range("F5").Value = 1.04
$ws.Range("G5").Value = 0.77
$ws.Range("D6").Value = 1.51
$ws.Range("G6").Value = 1.01
$ws.Range("E7").Value = 1.86
$ws.Range("F7").Value = 1.49
